$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of [cellRef, newValue] pairs describing every cell that changed,
# derived from the authoritative diff of the workbook XML.
$changes = @(
    @("D2","29.170.16"),
    @("E2","  +1.94%  "),
    @("D3","1.905.04"),
    @("E3","  +1.91%  "),
    @("E4","  -0.05%  "),
    @("D5","327.57"),
    @("E5","  +0.92%  "),
    @("D7","0.4651"),
    @("E7","  +0.77%  "),
    @("D8","0.3930"),
    @("E8","  +1.41%  "),
    @("B9","OKB"),
    @("C9","https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"),
    @("D9","46.86"),
    @("E9","  +1.09%  "),
    @("B10","Dogecoin"),
    @("C10","https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"),
    @("D10","0.07958"),
    @("E10","  +1.16%  "),
    @("B11","Polygon"),
    @("C11","https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"),
    @("D11","1.001"),
    @("E11","  +2.70%  "),
    @("B12","Solana"),
    @("C12","https://coinranking.com/coin/zNZHO_Sjf+solana-sol"),
    @("D12","22.35"),
    @("E12","  +1.91%  "),
    @("B13","WrappedEther"),
    @("C13","https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"),
    @("D13","1.962.72"),
    @("E13","  +4.91%  "),
    @("B14","Chainlink"),
    @("C14","https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"),
    @("D14","7.138"),
    @("E14","  +2.17%  "),
    @("B15","Polkadot"),
    @("C15","https://coinranking.com/coin/25W7FG7om+polkadot-dot"),
    @("D15","5.790"),
    @("E15","  +1.64%  "),
    @("B16","TRON"),
    @("C16","https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"),
    @("D16","0.06947"),
    @("E16","  -0.12%  "),
    @("B17","Litecoin"),
    @("C17","https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"),
    @("D17","88.77"),
    @("E17","  +0.71%  "),
    @("B18","BinanceUSD"),
    @("C18","https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"),
    @("D18","1.004"),
    @("E18","  -0.12%  "),
    @("B19","ShibaInu"),
    @("C19","https://coinranking.com/coin/xz24e0BjL+shibainu-shib"),
    @("D19","0.00001011"),
    @("E19","  +0.94%  "),
    @("B20","Avalanche"),
    @("C20","https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"),
    @("D20","17.16"),
    @("E20","  +2.15%  "),
    @("B21","Dai"),
    @("C21","https://coinranking.com/coin/MoTuySvg7+dai-dai"),
    @("D21","1.002"),
    @("E21","  -0.12%  "),
    @("B22","WrappedBTC"),
    @("C22","https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"),
    @("D22","29.191.00"),
    @("E22","  +2.03%  "),
    @("B23","Uniswap"),
    @("C23","https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"),
    @("D23","5.353"),
    @("E23","  +1.54%  "),
    @("B24","Cosmos"),
    @("C24","https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"),
    @("D24","11.07"),
    @("E24","  +0.46%  "),
    @("B25","WrappedliquidstakedEther2.0"),
    @("C25","https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"),
    @("D25","2.139.80"),
    @("E25","  +2.08%  "),
    @("B26","Toncoin"),
    @("C26","https://coinranking.com/coin/67YlI0K1b+toncoin-ton"),
    @("D26","2.062"),
    @("E26","  -2.32%  "),
    @("B27","Monero"),
    @("C27","https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"),
    @("D27","156.28"),
    @("E27","  +2.43%  "),
    @("B28","EthereumClassic"),
    @("C28","https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"),
    @("D28","19.52"),
    @("E28","  +1.47%  "),
    @("B29","InternetComputer(DFINITY)"),
    @("C29","https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"),
    @("D29","5.860"),
    @("E29","  +1.75%  "),
    @("B30","LidoDAOToken"),
    @("C30","https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"),
    @("D30","2.002"),
    @("E30","  +0.76%  "),
    @("B31","BitcoinCash"),
    @("C31","https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"),
    @("D31","119.42"),
    @("E31","  +0.25%  "),
    @("B32","Stellar"),
    @("C32","https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"),
    @("D32","0.09421"),
    @("E32","  +0.89%  "),
    @("B33","ImmutableX"),
    @("C33","https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"),
    @("D33","0.9234"),
    @("E33","  +0.66%  "),
    @("B34","Filecoin"),
    @("C34","https://coinranking.com/coin/ymQub4fuB+filecoin-fil"),
    @("D34","5.359"),
    @("E34","  +1.83%  "),
    @("B35","ARBITRUM"),
    @("C35","https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"),
    @("D35","1.345"),
    @("E35","  +0.75%  "),
    @("B36","HuobiToken"),
    @("C36","https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"),
    @("D36","3.257"),
    @("E36","  -2.09%  "),
    @("B37","Hedera"),
    @("C37","https://coinranking.com/coin/jad286TjB+hedera-hbar"),
    @("D37","0.05850"),
    @("E37","  +0.98%  "),
    @("B38","TrustWalletToken"),
    @("C38","https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"),
    @("D38","1.173"),
    @("E38","  +1.65%  "),
    @("B39","FraxShare"),
    @("C39","https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"),
    @("D39","7.998"),
    @("E39","  +3.41%  "),
    @("B40","VeChain"),
    @("C40","https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"),
    @("D40","0.02098"),
    @("E40","  +0.39%  "),
    @("B41","TheSandbox"),
    @("C41","https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"),
    @("D41","0.5753"),
    @("E41","  +2.22%  "),
    @("B42","Algorand"),
    @("C42","https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"),
    @("D42","0.1812"),
    @("E42","  +1.55%  "),
    @("B43","Aptos"),
    @("C43","https://coinranking.com/coin/HGYj5JCv5+aptos-apt"),
    @("D43","10.000"),
    @("E43","  +2.20%  "),
    @("B44","EnergySwap"),
    @("C44","https://coinranking.com/coin/SbWqqTui-+energyswap-ens"),
    @("D44","12.08"),
    @("E44","  +2.64%  "),
    @("B45","Decentraland"),
    @("C45","https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"),
    @("D45","0.5422"),
    @("E45","  +2.29%  "),
    @("B46","RenderToken"),
    @("C46","https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"),
    @("D46","2.225"),
    @("E46","  +2.82%  "),
    @("B47","Cronos"),
    @("C47","https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"),
    @("D47","0.07088"),
    @("E47","  -1.27%  "),
    @("B48","NEARProtocol"),
    @("C48","https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"),
    @("D48","1.880"),
    @("E48","  +2.64%  "),
    @("B49","MXToken"),
    @("C49","https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"),
    @("D49","2.545"),
    @("E49","  +5.73%  "),
    @("B50","Quant"),
    @("C50","https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"),
    @("D50","112.07"),
    @("E50","  -0.74%  "),
    @("B51","WEMIXToken"),
    @("C51","https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"),
    @("D51","1.071"),
    @("E51","  -6.23%  ")
)

foreach ($pair in $changes) {
    $ref = $pair[0]
    $val = $pair[1]
    $cell = $ws.Range($ref)
    $col = $ref.Substring(0,1)
    if ($col -eq "D" -or $col -eq "E") {
        # Columns D (Price) and E (Volume) hold numeric-looking text
        # (e.g. "10.000", "0.3930", "  +1.94%  "). Force the cell to be
        # treated as plain text so Excel does not auto-convert the
        # string into a number and strip formatting/precision.
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}
